$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.087.00"
$ws.Range("E2").Value = "  -2.60%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.665.49"
$ws.Range("E3").Value = "  -0.77%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "524.68"
$ws.Range("E5").Value = "  +0.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.56"
$ws.Range("E6").Value = "  -1.28%  "

$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.570"
$ws.Range("E8").Value = "  -1.01%  "

$ws.Range("E9").Value = "  +8.55%  "

$ws.Range("E10").Value = "  -2.41%  "

$ws.Range("E11").Value = "  -2.03%  "

$ws.Range("E12").Value = "  +1.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.135.70"
$ws.Range("E13").Value = "  -0.71%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "59.078.25"
$ws.Range("E14").Value = "  -2.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.08"
$ws.Range("E15").Value = "  -1.08%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000137"
$ws.Range("E16").Value = "  -1.81%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.669.97"
$ws.Range("E17").Value = "  -3.75%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "339.05"
$ws.Range("E18").Value = "  -3.45%  "

$ws.Range("E19").Value = "  -3.22%  "

$ws.Range("E20").Value = "  -2.29%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.42"
$ws.Range("E21").Value = "  +1.58%  "

$ws.Range("E22").Value = "  -0.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.39"
$ws.Range("E23").Value = "  +2.49%  "

$ws.Range("E24").Value = "  -0.91%  "

$ws.Range("E25").Value = "  -1.61%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  +0.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0803"
$ws.Range("E27").Value = "  -1.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.12"
$ws.Range("E28").Value = "  -2.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.68"
$ws.Range("E29").Value = "  -2.51%  "

$ws.Range("E30").Value = "  +0.05%  "

$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("E32").Value = "  -1.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.62"
$ws.Range("E33").Value = "  +1.56%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.16"
$ws.Range("E34").Value = "  -3.75%  "

$ws.Range("E35").Value = "  -4.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.897"
$ws.Range("E36").Value = "  -5.72%  "

$ws.Range("E37").Value = "  -0.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.94"
$ws.Range("E38").Value = "  +0.05%  "

$ws.Range("E39").Value = "  -5.81%  "

$ws.Range("E40").Value = "  -2.88%  "

$ws.Range("E41").Value = "  +0.79%  "

$ws.Range("E42").Value = "  +0.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "275.98"
$ws.Range("E43").Value = "  -2.16%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.91"
$ws.Range("E44").Value = "  -0.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0971"
$ws.Range("E45").Value = "  -1.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.67"
$ws.Range("E46").Value = "  +2.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.046.72"
$ws.Range("E47").Value = "  -4.12%  "

$ws.Range("E48").Value = "  -1.59%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.71"
$ws.Range("E49").Value = "  -3.43%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0229"
$ws.Range("E50").Value = "  -2.56%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.99"
$ws.Range("E51").Value = "  -1.47%  "

